# Apply the edits described by the commit diff:
#  - rename the shared string "proton" to "p" (cells I2:I10 on Sheet1)
#  - bold + center the header row (A1:K1)
#  - move the active selection to G18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename target "proton" -> "p" for all data rows (I2:I10)
$ws.Range("I2:I10").Value = "p"

# Header row (A1:K1): make it bold while keeping the existing centered alignment
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# Move the selection to match the saved workbook state
$ws.Range("G18").Select()
